$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.2633423308076207
$ws.Range("D2").Value = 0.7947374909435916

$ws.Range("C3").Value = 0.1596128116582426
$ws.Range("D3").Value = 0.8746424273097646

$ws.Range("C4").Value = 1.963336204694945
$ws.Range("D4").Value = 0.06237104623854806
$ws.Range("G4").Value = "No"

$ws.Range("C5").Value = 2.239979968767079
$ws.Range("D5").Value = 0.0355207219835183

$ws.Range("C6").Value = 0.3708920834244382
$ws.Range("D6").Value = 0.7142653835892374

$ws.Range("C7").Value = 2.451826326039056
$ws.Range("D7").Value = 0.02261802533419632

$ws.Range("C8").Value = 2.249545626871841
$ws.Range("D8").Value = 0.03481651128637697

$ws.Range("C9").Value = 2.321485767588634
$ws.Range("D9").Value = 0.02991569475022238

$ws.Range("C10").Value = 3.486425884470581
$ws.Range("D10").Value = 0.002091076247400681

$ws.Range("C11").Value = -0.2182745357496441
$ws.Range("D11").Value = 0.829227234442641
